# actualiza datos provincias 2022-05-03
# Updates the consolidated COVID-19 Spain dataset rows 823-855 (2022-04-29 .. 2022-05-10)
# with refreshed province-level figures (hospitalized, deceased, daily_deaths, etc.)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(823, 6).Value = 3994
$ws.Cells.Item(823, 10).Value = 141660
$ws.Cells.Item(823, 22).Value = 48717
$ws.Cells.Item(823, 24).Value = 8029.7
$ws.Cells.Item(823, 25).Value = 3458.3
$ws.Cells.Item(823, 26).Value = 24208
$ws.Cells.Item(824, 6).Value = 2969
$ws.Cells.Item(824, 10).Value = 100563
$ws.Cells.Item(824, 22).Value = -41097
$ws.Cells.Item(824, 23).Value = -29
$ws.Cells.Item(824, 24).Value = -5139
$ws.Cells.Item(824, 25).Value = -2416
$ws.Cells.Item(824, 26).Value = -16912
$ws.Cells.Item(825, 10).Value = 100592
$ws.Cells.Item(825, 24).Value = 2549.7
$ws.Cells.Item(825, 25).Value = -2414.3
$ws.Cells.Item(825, 26).Value = -16900
$ws.Cells.Item(826, 10).Value = 113949
$ws.Cells.Item(826, 23).Value = 13.3
$ws.Cells.Item(826, 24).Value = -9237
$ws.Cells.Item(826, 25).Value = -508.7
$ws.Cells.Item(826, 26).Value = -3561
$ws.Cells.Item(827, 10).Value = 114025
$ws.Cells.Item(827, 25).Value = -506.6
$ws.Cells.Item(827, 26).Value = -3546
$ws.Cells.Item(828, 10).Value = 100669
$ws.Cells.Item(828, 23).Value = -11.7
$ws.Cells.Item(828, 25).Value = -2187.3
$ws.Cells.Item(828, 26).Value = -15311
$ws.Cells.Item(829, 10).Value = 114180
$ws.Cells.Item(829, 23).Value = 13.4
$ws.Cells.Item(829, 25).Value = 3033.9
$ws.Cells.Item(829, 26).Value = 21237
$ws.Cells.Item(830, 10).Value = 114255
$ws.Cells.Item(830, 25).Value = -3915
$ws.Cells.Item(830, 26).Value = -27405
$ws.Cells.Item(831, 10).Value = 100864
$ws.Cells.Item(831, 23).Value = -11.7
$ws.Cells.Item(832, 10).Value = 100905
$ws.Cells.Item(833, 10).Value = 114336
$ws.Cells.Item(833, 23).Value = 13.3
$ws.Cells.Item(834, 10).Value = 114495
$ws.Cells.Item(834, 22).Value = 159
$ws.Cells.Item(834, 24).Value = 4543.7
$ws.Cells.Item(834, 25).Value = 67.1
$ws.Cells.Item(834, 26).Value = 470
$ws.Cells.Item(835, 10).Value = 101079
$ws.Cells.Item(835, 23).Value = -11.7
$ws.Cells.Item(835, 24).Value = 58
$ws.Cells.Item(835, 25).Value = 58.6
$ws.Cells.Item(835, 26).Value = 410
$ws.Cells.Item(836, 10).Value = 101095
$ws.Cells.Item(836, 24).Value = -4413.7
$ws.Cells.Item(836, 25).Value = -1869.3
$ws.Cells.Item(836, 26).Value = -13085
$ws.Cells.Item(837, 10).Value = 101121
$ws.Cells.Item(837, 22).Value = 26
$ws.Cells.Item(837, 24).Value = -4458
$ws.Cells.Item(837, 25).Value = -1876.3
$ws.Cells.Item(837, 26).Value = -13134
$ws.Cells.Item(838, 10).Value = 101148
$ws.Cells.Item(838, 22).Value = 27
$ws.Cells.Item(838, 24).Value = 23
$ws.Cells.Item(838, 25).Value = 40.6
$ws.Cells.Item(838, 26).Value = 284
$ws.Cells.Item(839, 10).Value = 101192
$ws.Cells.Item(839, 22).Value = 44
$ws.Cells.Item(839, 24).Value = 32.3
$ws.Cells.Item(839, 25).Value = 41
$ws.Cells.Item(839, 26).Value = 287
$ws.Cells.Item(840, 10).Value = 101217
$ws.Cells.Item(840, 22).Value = 25
$ws.Cells.Item(840, 24).Value = 32
$ws.Cells.Item(840, 25).Value = -1874.1
$ws.Cells.Item(840, 26).Value = -13119
$ws.Cells.Item(841, 10).Value = 114755
$ws.Cells.Item(841, 23).Value = 13.4
$ws.Cells.Item(841, 24).Value = 4535.7
$ws.Cells.Item(841, 25).Value = 37.1
$ws.Cells.Item(841, 26).Value = 260
$ws.Cells.Item(842, 10).Value = 92073
$ws.Cells.Item(842, 23).Value = -19.8
$ws.Cells.Item(842, 24).Value = -3039.7
$ws.Cells.Item(842, 25).Value = -1286.6
$ws.Cells.Item(842, 26).Value = -9006
$ws.Cells.Item(843, 10).Value = 105737
$ws.Cells.Item(843, 22).Value = 13664
$ws.Cells.Item(843, 23).Value = 14.8
$ws.Cells.Item(843, 24).Value = 1506.7
$ws.Cells.Item(843, 25).Value = 663.1
$ws.Cells.Item(843, 26).Value = 4642
$ws.Cells.Item(844, 10).Value = 115048
$ws.Cells.Item(844, 22).Value = 9311
$ws.Cells.Item(844, 23).Value = 8.8
$ws.Cells.Item(844, 24).Value = 97.7
$ws.Cells.Item(844, 25).Value = 1989.6
$ws.Cells.Item(844, 26).Value = 13927
$ws.Cells.Item(845, 10).Value = 101535
$ws.Cells.Item(845, 22).Value = -13513
$ws.Cells.Item(845, 23).Value = -11.7
$ws.Cells.Item(845, 24).Value = 3154
$ws.Cells.Item(845, 25).Value = 55.3
$ws.Cells.Item(845, 26).Value = 387
$ws.Cells.Item(846, 10).Value = 101587
$ws.Cells.Item(846, 22).Value = 52
$ws.Cells.Item(846, 23).Value = 0.1
$ws.Cells.Item(846, 24).Value = -1383.3
$ws.Cells.Item(846, 25).Value = 56.4
$ws.Cells.Item(846, 26).Value = 395
$ws.Cells.Item(847, 2).Value = 3832
$ws.Cells.Item(847, 4).Value = 6169077
$ws.Cells.Item(847, 6).Value = 4839
$ws.Cells.Item(847, 7).Value = 127721
$ws.Cells.Item(847, 9).Value = 202
$ws.Cells.Item(847, 10).Value = 115151
$ws.Cells.Item(847, 11).Value = 8977333
$ws.Cells.Item(847, 14).Value = 133782
$ws.Cells.Item(847, 15).Value = 1510647
$ws.Cells.Item(847, 18).Value = 1451381
$ws.Cells.Item(847, 19).Value = 215806.7
$ws.Cells.Item(847, 22).Value = 13564
$ws.Cells.Item(847, 23).Value = 13.4
$ws.Cells.Item(847, 24).Value = 34.3
$ws.Cells.Item(847, 25).Value = 1990.6
$ws.Cells.Item(847, 26).Value = 13934
$ws.Cells.Item(848, 2).Value = 10755
$ws.Cells.Item(848, 4).Value = 6045840
$ws.Cells.Item(848, 6).Value = 4953
$ws.Cells.Item(848, 7).Value = 415
$ws.Cells.Item(848, 9).Value = 195
$ws.Cells.Item(848, 10).Value = 115233
$ws.Cells.Item(848, 11).Value = 8999266
$ws.Cells.Item(848, 14).Value = 139781
$ws.Cells.Item(848, 15).Value = 80125
$ws.Cells.Item(848, 18).Value = 21933
$ws.Cells.Item(848, 19).Value = 11446.4
$ws.Cells.Item(848, 22).Value = 82
$ws.Cells.Item(848, 24).Value = 4566
$ws.Cells.Item(848, 25).Value = 68.3
$ws.Cells.Item(848, 26).Value = 478
$ws.Cells.Item(849, 2).Value = 4311
$ws.Cells.Item(849, 4).Value = 6177655
$ws.Cells.Item(849, 6).Value = 4306
$ws.Cells.Item(849, 7).Value = 127953
$ws.Cells.Item(849, 9).Value = 160
$ws.Cells.Item(849, 10).Value = 101712
$ws.Cells.Item(849, 11).Value = 7561449
$ws.Cells.Item(849, 14).Value = 122288
$ws.Cells.Item(849, 15).Value = 1442293
$ws.Cells.Item(849, 18).Value = -1437817
$ws.Cells.Item(849, 19).Value = 206041.9
$ws.Cells.Item(849, 22).Value = -13521
$ws.Cells.Item(849, 23).Value = -11.7
$ws.Cells.Item(849, 24).Value = 41.7
$ws.Cells.Item(849, 25).Value = 1377
$ws.Cells.Item(849, 26).Value = 9639
$ws.Cells.Item(850, 2).Value = 4170
$ws.Cells.Item(850, 4).Value = 6180726
$ws.Cells.Item(850, 6).Value = 4523
$ws.Cells.Item(850, 7).Value = 127599
$ws.Cells.Item(850, 9).Value = 166
$ws.Cells.Item(850, 10).Value = 112730
$ws.Cells.Item(850, 11).Value = 8798609
$ws.Cells.Item(850, 14).Value = 1355635
$ws.Cells.Item(850, 15).Value = 1237493
$ws.Cells.Item(850, 18).Value = 1237160
$ws.Cells.Item(850, 19).Value = 176784.7
$ws.Cells.Item(850, 23).Value = 10.8
$ws.Cells.Item(850, 24).Value = -807
$ws.Cells.Item(850, 25).Value = 999
$ws.Cells.Item(850, 26).Value = 6993
$ws.Cells.Item(851, 2).Value = 8006
$ws.Cells.Item(851, 4).Value = 6056803
$ws.Cells.Item(851, 6).Value = 4439
$ws.Cells.Item(851, 7).Value = 40
$ws.Cells.Item(851, 9).Value = 173
$ws.Cells.Item(851, 10).Value = 112818
$ws.Cells.Item(851, 11).Value = 8821275
$ws.Cells.Item(851, 14).Value = 1375926
$ws.Cells.Item(851, 15).Value = -137993
$ws.Cells.Item(851, 18).Value = 22666
$ws.Cells.Item(851, 19).Value = -19713.3
$ws.Cells.Item(851, 24).Value = -805
$ws.Cells.Item(851, 25).Value = -318.6
$ws.Cells.Item(851, 26).Value = -2230
$ws.Cells.Item(852, 2).Value = 2485
$ws.Cells.Item(852, 4).Value = 6059278
$ws.Cells.Item(852, 6).Value = 3886
$ws.Cells.Item(852, 7).Value = 53
$ws.Cells.Item(852, 9).Value = 99
$ws.Cells.Item(852, 10).Value = 99209
$ws.Cells.Item(852, 11).Value = 7366184
$ws.Cells.Item(852, 14).Value = -86481
$ws.Cells.Item(852, 15).Value = -152704
$ws.Cells.Item(852, 18).Value = -1455091
$ws.Cells.Item(852, 19).Value = -21814.9
$ws.Cells.Item(852, 23).Value = -12.1
$ws.Cells.Item(852, 25).Value = -332.3
$ws.Cells.Item(852, 26).Value = -2326
$ws.Cells.Item(853, 2).Value = 1156
$ws.Cells.Item(853, 4).Value = 6189063
$ws.Cells.Item(853, 6).Value = 3821
$ws.Cells.Item(853, 7).Value = 127646
$ws.Cells.Item(853, 9).Value = 144
$ws.Cells.Item(853, 10).Value = 71432
$ws.Cells.Item(853, 11).Value = 7375133
$ws.Cells.Item(853, 14).Value = -86046
$ws.Cells.Item(853, 15).Value = -150819
$ws.Cells.Item(853, 18).Value = 8949
$ws.Cells.Item(853, 19).Value = -21545.6
$ws.Cells.Item(853, 22).Value = -27777
$ws.Cells.Item(853, 23).Value = -28
$ws.Cells.Item(853, 24).Value = -13766
$ws.Cells.Item(853, 25).Value = -4307.9
$ws.Cells.Item(853, 26).Value = -30155
$ws.Cells.Item(854, 15).Value = -844192
$ws.Cells.Item(854, 18).Value = 758008
$ws.Cells.Item(854, 19).Value = -120598.9
$ws.Cells.Item(854, 22).Value = 6961
$ws.Cells.Item(854, 23).Value = 9.7
$ws.Cells.Item(854, 24).Value = -11475
$ws.Cells.Item(854, 25).Value = -5251.1
$ws.Cells.Item(854, 26).Value = -36758
$ws.Cells.Item(855, 15).Value = -4327403
$ws.Cells.Item(855, 19).Value = -618200.4
$ws.Cells.Item(855, 24).Value = -17790
$ws.Cells.Item(855, 25).Value = -9913.4
$ws.Cells.Item(855, 26).Value = -69394
